# Remove the two placeholder "# Dummy Row (Ignored by ExcelBinder)" rows
# from the "Skill" sheet (rows 2 and 3), shifting the real data up so it
# now occupies rows 2 and 3. Then select cell E3 on the Skill sheet and
# activate it, making it the active/selected tab (previously "Item" was
# the selected tab).

$wb = $excel.ActiveWorkbook
$skill = $wb.Worksheets.Item("Skill")

# Delete rows 2:3 (the dummy placeholder rows), real data shifts up.
$skill.Rows("2:3").Delete()

# Activate the Skill sheet and select E3, matching the new active tab /
# selection recorded in the sheet view.
$skill.Activate()
$skill.Range("E3").Select()
